# Update return shipments data (auto)
#
# The upstream AfterShip export was refreshed: several return-shipment rows
# were re-ordered/refreshed with newer checkpoint timestamps, one brand-new
# shipment (Direct Freight Express, AU) was inserted at row 31 pushing the
# remaining "csv_importer" rows down by one, and four stale
# spanish-seur-ftp / swiss-post rows that no longer appear in the refreshed
# export were dropped from the bottom of the sheet (old rows 36-39),
# shrinking the used range from A1:M39 to A1:M35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only cells whose value actually changes are written. Purely-numeric-looking
# values (tracking numbers, order ids, sales office ids, ...) are written with
# a leading apostrophe to force text storage so they keep their leading zeros /
# full integer precision instead of being coerced to numbers.

# Row 11
$ws.Range("A11").Value = "00KLOK5U"
$ws.Range("E11").Value = "'0031009885087"
$ws.Range("I11").Value = "2026-02-12T16:21:03+01:00"
$ws.Range("K11").Value = "2026-02-12T15:50:46+00:00"
$ws.Range("L11").Value = "00KLOK5U"
$ws.Range("M11").Value = "{`"external_order_id`": `"0031009885087`", `"sales_office_id`": `"0303`"}"

# Row 12
$ws.Range("A12").Value = "00KLOJSZ"
$ws.Range("E12").Value = "'0031009884566"
$ws.Range("I12").Value = "2026-02-12T11:36:34+01:00"
$ws.Range("K12").Value = "2026-02-12T13:50:37+00:00"
$ws.Range("L12").Value = "00KLOJSZ"
$ws.Range("M12").Value = "{`"external_order_id`": `"0031009884566`", `"sales_office_id`": `"0303`"}"

# Row 13
$ws.Range("A13").Value = "00KLOK0B"
$ws.Range("K13").Value = "2026-02-12T13:50:38+00:00"
$ws.Range("L13").Value = "00KLOK0B"

# Row 16
$ws.Range("A16").Value = "ZI6799YB"
$ws.Range("L16").Value = "ZI6799YB"

# Row 17
$ws.Range("A17").Value = "ZI679A26"
$ws.Range("L17").Value = "ZI679A26"

# Row 19
$ws.Range("A19").Value = "ZI6799VW"
$ws.Range("E19").Value = "'0032005538280"
$ws.Range("I19").Value = "2026-02-12T15:09:56+01:00"
$ws.Range("K19").Value = "2026-02-12T14:23:45+00:00"
$ws.Range("L19").Value = "ZI6799VW"
$ws.Range("M19").Value = "{`"external_order_id`": `"0032005538280`", `"sales_office_id`": `"0303`"}"

# Row 20
$ws.Range("A20").Value = "ZI679A2O"
$ws.Range("E20").Value = "'0032005538512"
$ws.Range("I20").Value = "2026-02-12T12:37:12+01:00"
$ws.Range("K20").Value = "2026-02-12T13:47:11+00:00"
$ws.Range("L20").Value = "ZI679A2O"
$ws.Range("M20").Value = "{`"external_order_id`": `"0032005538512`", `"sales_office_id`": `"0303`"}"

# Row 22
$ws.Range("A22").Value = "'81225942281"
$ws.Range("I22").Value = "2026-02-11T12:20:19+02:00"
$ws.Range("K22").Value = "2026-02-12T12:20:38+00:00"
$ws.Range("L22").Value = "'81225942281"

# Row 23
$ws.Range("A23").Value = "'81225941052"
$ws.Range("I23").Value = "2026-02-11T12:19:10+02:00"
$ws.Range("K23").Value = "2026-02-12T12:20:33+00:00"
$ws.Range("L23").Value = "'81225941052"

# Row 31
$ws.Range("A31").Value = "'3366515579266"
$ws.Range("B31").Value = "directfreight-au-ref"
$ws.Range("C31").Value = "Direct Freight Express"
$ws.Range("E31").Value = "250-3604430-4903019"
$ws.Range("F31").Value = "'4272"
$ws.Range("G31").Value = "api"
$ws.Range("I31").Value = "2026-02-13T07:51:00+11:00"
$ws.Range("J31").Value = "Warrnambool, Australia"
$ws.Range("K31").Value = "2026-02-12T21:02:10+00:00"
$ws.Range("L31").Value = "'3366515579266"
$ws.Range("M31").Value = "{`"external_order_id`": `"250-3604430-4903019`", `"sales_office_id`": `"4272`"}"

# Row 32
$ws.Range("A32").Value = "1Z0JA1729022187602"
$ws.Range("B32").Value = "ups-api"
$ws.Range("C32").Value = "ups-api"
$ws.Range("E32").Value = "'6001609505"
$ws.Range("I32").Value = "2026-02-11T11:58:21-08:00"
$ws.Range("J32").Value = "REDLANDS, CA, 92374, US, United States"
$ws.Range("K32").Value = "2026-02-12T12:14:08+00:00"
$ws.Range("L32").Value = "1Z0JA1729022187602"
$ws.Range("M32").Value = "{`"courier`": `"UPS`", `"custom_1`": `"returns_ups`"}"

# Row 33
$ws.Range("A33").Value = "'996011434800723580"
$ws.Range("B33").Value = "swiss-post"
$ws.Range("C33").Value = "Swiss Post"
$ws.Range("E33").Value = "'6001610445"
$ws.Range("I33").Value = "2026-02-12T10:27:29+01:00"
$ws.Range("J33").Value = "CH, Switzerland"
$ws.Range("K33").Value = "2026-02-12T10:14:17+00:00"
$ws.Range("L33").Value = "'996011434800723580"
$ws.Range("M33").Value = "{`"courier`": `"swiss-post`"}"

# Row 35
$ws.Range("A35").Value = "'996011434800724656"
$ws.Range("B35").Value = "swiss-post"
$ws.Range("C35").Value = "Swiss Post"
$ws.Range("E35").Value = "'6001609189"
$ws.Range("I35").Value = "2026-02-12T09:53:31+01:00"
$ws.Range("J35").Value = "CH, Switzerland"
$ws.Range("K35").Value = "2026-02-12T09:14:18+00:00"
$ws.Range("L35").Value = "'996011434800724656"
$ws.Range("M35").Value = "{`"courier`": `"swiss-post`"}"

# Remove the now-obsolete trailing rows (old rows 36-39 merged/removed)
$ws.Rows("36:39").Delete()
